$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Impl1" data row (row 2: values -3.65 / 291201.116637 / 58830.508095).
# This shifts every row below it up by one (carrying formatting/styles along),
# so what used to be row 3 becomes row 2, row 4 becomes row 3, row 5 becomes
# row 4 and row 6 becomes row 5.
$ws.Rows.Item(2).Delete()

# After the shift the "Name" labels in column A need to slide up one slot too
# (each row now shows the label that used to belong to the implementation one
# step earlier), with the final row's label renamed to "Impl4".
$ws.Range("A2").Value = "impl1"
$ws.Range("A3").Value = "impl2"
$ws.Range("A4").Value = "impl3"
$ws.Range("A5").Value = "Impl4"

# Reposition the selection (cursor) as in the authored change.
[void]$ws.Range("I13").Select()
